# Fruta / hortaliza, semanal
# Insert a new week's worth of price data (4 rows) for Melón at
# Terminal Hortofrutícola Agro Chillán, pushing the existing history down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 245 (shifts old rows 245:261 down to 249:265).
$ws.Range("A245:R248").EntireRow.Insert()

# Values shared by every row in this data block.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$catId     = 100112027
$categoria = "Melón"
$unidadCom = "$/unidad"
$kgUnid    = 1
$clasif    = "Hortaliza"
$fecha     = 44931

function Set-Fila($row, $variedad, $calidad, $volumen, $min, $max, $prom, $origen) {
    $ws.Range("A$row").Value = $mercadoId
    $ws.Range("B$row").Value = $mercado
    $ws.Range("C$row").Value = $region
    $ws.Range("D$row").Value = $fecha
    $ws.Range("E$row").Value = $codreg
    $ws.Range("F$row").Value = $catId
    $ws.Range("G$row").Value = $categoria
    $ws.Range("H$row").Value = $variedad
    $ws.Range("I$row").Value = $calidad
    $ws.Range("J$row").Value = $volumen
    $ws.Range("K$row").Value = $min
    $ws.Range("L$row").Value = $max
    $ws.Range("M$row").Value = $prom
    $ws.Range("N$row").Value = $unidadCom
    $ws.Range("O$row").Value = $origen
    $ws.Range("P$row").Value = $prom
    $ws.Range("Q$row").Value = $kgUnid
    $ws.Range("R$row").Value = $clasif
}

Set-Fila 245 "Calameño" "Extra"   500 1000 1000 1000 "Región de O'Higgins"
Set-Fila 246 "Calameño" "Primera" 500 800  800  800  "Región de O'Higgins"
Set-Fila 247 "Tuna"     "Extra"   500 1000 1000 1000 "Región de O'Higgins"
Set-Fila 248 "Tuna"     "Primera" 500 800  800  800  "Región de O'Higgins"
